$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 421.7027586666666
    "H2" = 1265.108276
    "I2" = 0.8230205889389596
    "J2" = 0.8230205889389596
    "M2" = 211.2725676666666
    "N2" = 633.8177029999999
    "O2" = 0.6324644927232657
    "P2" = 0.6324644927232657
    "Q2" = 89094.22461562332
    "R2" = 801848.0215406099
    "S2" = 0.5205312992840824
    "T2" = 0.5205312992840824
    "G3" = 421.7027586666666
    "H3" = 1265.108276
    "I3" = 0.8230205889389596
    "J3" = 0.8230205889389596
    "M3" = 59.36675400000001
    "N3" = 178.100262
    "O3" = 0.1777200152765546
    "P3" = 0.1777200152765546
    "Q3" = 25035.12393488537
    "R3" = 225316.1154139683
    "S3" = 0.1462672316391508
    "T3" = 0.1462672316391508
    "G4" = 421.7027586666666
    "H4" = 1265.108276
    "I4" = 0.8230205889389596
    "J4" = 0.8230205889389596
    "M4" = 0.4593846666666666
    "N4" = 1.378154
    "O4" = 0.001375211620595172
    "P4" = 0.001375211620595172
    "Q4" = 193.7237812225004
    "R4" = 1743.514031002504
    "S4" = 0.00113182747789794
    "T4" = 0.00113182747789794
    "G5" = 421.7027586666666
    "H5" = 1265.108276
    "I5" = 0.8230205889389596
    "J5" = 0.8230205889389596
    "M5" = 62.94782133333333
    "N5" = 188.843464
    "O5" = 0.1884402803795846
    "P5" = 0.1884402803795846
    "Q5" = 26545.26990832311
    "R5" = 238907.429174908
    "S5" = 0.1550902305378284
    "T5" = 0.1550902305378284
    "I6" = 0.136040387754698
    "J6" = 0.136040387754698
    "M6" = 211.2725676666666
    "N6" = 633.8177029999999
    "O6" = 0.6324644927232657
    "P6" = 0.6324644927232657
    "Q6" = 14726.7432021831
    "R6" = 132540.6888196479
    "S6" = 0.08604071483115144
    "T6" = 0.08604071483115144
    "I7" = 0.136040387754698
    "J7" = 0.136040387754698
    "M7" = 59.36675400000001
    "N7" = 178.100262
    "O7" = 0.1777200152765546
    "P7" = 0.1777200152765546
    "Q7" = 4138.156460920956
    "R7" = 37243.40814828861
    "S7" = 0.02417709978999334
    "T7" = 0.02417709978999334
    "I8" = 0.136040387754698
    "J8" = 0.136040387754698
    "M8" = 0.4593846666666666
    "N8" = 1.378154
    "O8" = 0.001375211620595172
    "P8" = 0.001375211620595172
    "Q8" = 32.02138399574088
    "R8" = 288.192455961668
    "S8" = 0.0001870843221105339
    "T8" = 0.0001870843221105339
    "I9" = 0.136040387754698
    "J9" = 0.136040387754698
    "M9" = 62.94782133333333
    "N9" = 188.843464
    "O9" = 0.1884402803795846
    "P9" = 0.1884402803795846
    "Q9" = 4387.774570788077
    "R9" = 39489.97113709269
    "S9" = 0.02563548881144271
    "T9" = 0.0256354888114427
    "G10" = 1.088159666666667
    "H10" = 3.264479
    "I10" = 0.002123718167154624
    "J10" = 0.002123718167154624
    "M10" = 211.2725676666666
    "N10" = 633.8177029999999
    "O10" = 0.6324644927232657
    "P10" = 0.6324644927232657
    "Q10" = 229.8982868079707
    "R10" = 2069.084581271737
    "S10" = 0.001343176333276633
    "T10" = 0.001343176333276633
    "G11" = 1.088159666666667
    "H11" = 3.264479
    "I11" = 0.002123718167154624
    "J11" = 0.002123718167154624
    "M11" = 59.36675400000001
    "N11" = 178.100262
    "O11" = 0.1777200152765546
    "P11" = 0.1777200152765546
    "Q11" = 64.600507243722
    "R11" = 581.4045651934981
    "S11" = 0.0003774272251098163
    "T11" = 0.0003774272251098163
    "G12" = 1.088159666666667
    "H12" = 3.264479
    "I12" = 0.002123718167154624
    "J12" = 0.002123718167154624
    "M12" = 0.4593846666666666
    "N12" = 1.378154
    "O12" = 0.001375211620595172
    "P12" = 0.001375211620595172
    "Q12" = 0.4998838657517777
    "R12" = 4.498954791766
    "S12" = 0.000002920561902340119
    "T12" = 0.000002920561902340119
    "G13" = 1.088159666666667
    "H13" = 3.264479
    "I13" = 0.002123718167154624
    "J13" = 0.002123718167154624
    "M13" = 62.94782133333333
    "N13" = 188.843464
    "O13" = 0.1884402803795846
    "P13" = 0.1884402803795846
    "Q13" = 68.49728027947289
    "R13" = 616.475522515256
    "S13" = 0.0004001940468658349
    "T13" = 0.0004001940468658349
    "G14" = 19.43698366666667
    "H14" = 58.310951
    "I14" = 0.03793439197579861
    "J14" = 0.03793439197579861
    "M14" = 211.2725676666666
    "N14" = 633.8177029999999
    "O14" = 0.6324644927232657
    "P14" = 0.6324644927232657
    "Q14" = 4106.501446951727
    "R14" = 36958.51302256555
    "S14" = 0.02399215597773899
    "T14" = 0.02399215597773899
    "G15" = 19.43698366666667
    "H15" = 58.310951
    "I15" = 0.03793439197579861
    "J15" = 0.03793439197579861
    "M15" = 59.36675400000001
    "N15" = 178.100262
    "O15" = 0.1777200152765546
    "P15" = 0.1777200152765546
    "Q15" = 1153.910627841018
    "R15" = 10385.19565056916
    "S15" = 0.00674170072144574
    "T15" = 0.006741700721445739
    "G16" = 19.43698366666667
    "H16" = 58.310951
    "I16" = 0.03793439197579861
    "J16" = 0.03793439197579861
    "M16" = 0.4593846666666666
    "N16" = 1.378154
    "O16" = 0.001375211620595172
    "P16" = 0.001375211620595172
    "Q16" = 8.92905226271711
    "R16" = 80.361470364454
    "S16" = 0.0000521678166653305
    "T16" = 0.0000521678166653305
    "G17" = 19.43698366666667
    "H17" = 58.310951
    "I17" = 0.03793439197579861
    "J17" = 0.03793439197579861
    "M17" = 62.94782133333333
    "N17" = 188.843464
    "O17" = 0.1884402803795846
    "P17" = 0.1884402803795846
    "Q17" = 1223.515775108252
    "R17" = 11011.64197597426
    "S17" = 0.007148367459948556
    "T17" = 0.007148367459948555
    "G18" = 0.451366
    "H18" = 1.354098
    "I18" = 0.0008809131633892397
    "J18" = 0.0008809131633892397
    "M18" = 211.2725676666666
    "N18" = 633.8177029999999
    "O18" = 0.6324644927232657
    "P18" = 0.6324644927232657
    "Q18" = 95.36125377743265
    "R18" = 858.2512839968939
    "S18" = 0.0005571462970162227
    "T18" = 0.0005571462970162227
    "G19" = 0.451366
    "H19" = 1.354098
    "I19" = 0.0008809131633892397
    "J19" = 0.0008809131633892397
    "M19" = 59.36675400000001
    "N19" = 178.100262
    "O19" = 0.1777200152765546
    "P19" = 0.1777200152765546
    "Q19" = 26.796134285964
    "R19" = 241.165208573676
    "S19" = 0.0001565559008548537
    "T19" = 0.0001565559008548537
    "G20" = 0.451366
    "H20" = 1.354098
    "I20" = 0.0008809131633892397
    "J20" = 0.0008809131633892397
    "M20" = 0.4593846666666666
    "N20" = 1.378154
    "O20" = 0.001375211620595172
    "P20" = 0.001375211620595172
    "Q20" = 0.2073506194546666
    "R20" = 1.866155575092
    "S20" = 0.000001211442019028136
    "T20" = 0.000001211442019028136
    "G21" = 0.451366
    "H21" = 1.354098
    "I21" = 0.0008809131633892397
    "J21" = 0.0008809131633892397
    "M21" = 62.94782133333333
    "N21" = 188.843464
    "O21" = 0.1884402803795846
    "P21" = 0.1884402803795846
    "Q21" = 28.41250632394133
    "R21" = 255.712556915472
    "S21" = 0.0001659995234991352
    "T21" = 0.0001659995234991352
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Host "Updated $($updates.Count) cells"